$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "G" (HACİZ GÜNÜ) column so the
# old G/H columns shift right to H/I, then populate the new column with
# the "ŞEHİR" (city) classification for each row.
$ws.Columns("G").Insert()

$ws.Range("G1").Value = "ŞEHİR"
$ws.Range("G2").Value = "İSTANBUL"
$ws.Range("G3").Value = "ANKARA"
$ws.Range("G4").Value = "ŞANLIURFA"
$ws.Range("G5").Value = "ŞANLIURFA"
$ws.Range("G6").Value = "ANKARA"
$ws.Range("G7").Value = "DENİZLİ"
$ws.Range("G8").Value = "İZMİR"
$ws.Range("G9").Value = "İZMİR"
$ws.Range("G10").Value = "BURSA"
$ws.Range("G11").Value = "BALIKESİR"

$ws.Columns("G").ColumnWidth = 8.5

$ws.Range("I12").Select()
